$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared string "PWM输出" -> "PWM输出检查口" in place (G13 is the sole
# cell using that string, so touching it first keeps it at the same
# sharedStrings index instead of minting a brand-new one at the end).
$ws.Range("G13").Value = "PWM输出检查口"

# New values for the previously-empty G column cells.
$ws.Range("G12").Value = "ADC采集检查口"
$ws.Range("G1").Value = "临时功能"

# Column G needs an explicit width (splits the old 5:16384 default-width
# run into 5:6, 7 (custom), 8:16384).
$ws.Range("G1").EntireColumn.ColumnWidth = 13

# View: zoom 190% -> 175%, and move the active selection from D18 to G2.
$excel.ActiveWindow.Zoom = 175
$ws.Range("G2").Select()
